# Append new daily COVID death rows (234-273) to Sheet1, mirroring the
# upstream OpenData CSV refresh: dates 2021-06-06 .. 2021-07-15 (serial
# 44353..44392) with DeathCovid / DeathWithCovid / Total cumulative counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Date (serial), B=DeathCovid, C=DeathWithCovid, D=Total
$newRows = @(
    @(44353, 12414, 2426, 14840),
    @(44354, 12423, 2426, 14849),
    @(44355, 12430, 2426, 14856),
    @(44356, 12433, 2426, 14859),
    @(44357, 12433, 2426, 14859),
    @(44358, 12436, 2427, 14863),
    @(44359, 12439, 2427, 14866),
    @(44360, 12441, 2428, 14869),
    @(44361, 12446, 2430, 14876),
    @(44362, 12456, 2433, 14889),
    @(44363, 12464, 2437, 14901),
    @(44364, 12478, 2442, 14920),
    @(44365, 12486, 2446, 14932),
    @(44366, 12492, 2450, 14942),
    @(44367, 12496, 2450, 14946),
    @(44368, 12502, 2450, 14952),
    @(44369, 12502, 2450, 14952),
    @(44370, 12502, 2450, 14952),
    @(44371, 12505, 2451, 14956),
    @(44372, 12505, 2451, 14956),
    @(44373, 12505, 2451, 14956),
    @(44374, 12509, 2451, 14960),
    @(44375, 12510, 2451, 14961),
    @(44376, 12510, 2451, 14961),
    @(44377, 12511, 2451, 14962),
    @(44378, 12511, 2451, 14962),
    @(44379, 12513, 2451, 14964),
    @(44380, 12513, 2451, 14964),
    @(44381, 12513, 2451, 14964),
    @(44382, 12514, 2451, 14965),
    @(44383, 12516, 2451, 14967),
    @(44384, 12516, 2451, 14967),
    @(44385, 12517, 2451, 14968),
    @(44386, 12519, 2452, 14971),
    @(44387, 12519, 2452, 14971),
    @(44388, 12521, 2452, 14973),
    @(44389, 12522, 2453, 14975),
    @(44390, 12523, 2453, 14976),
    @(44391, 12524, 2454, 14978),
    @(44392, 12524, 2454, 14978)
)

$startRow = 234
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $vals[0]
    $ws.Cells.Item($r, 2).Value2 = $vals[1]
    $ws.Cells.Item($r, 3).Value2 = $vals[2]
    $ws.Cells.Item($r, 4).Value2 = $vals[3]
}

$lastRow = $startRow + $newRows.Count - 1
$ws.Range("A" + $lastRow).Select()
